# Scheduled market-data refresh: update currentAveragePrice* / Leve*Price* /
# LeveProfit* columns (H:N) on affected Leve rows across the per-job sheets,
# reflecting the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 2900
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2900
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2900
$ws.Range("N16").Value = -3360
$ws.Range("M16").ClearContents()

$ws.Range("H17").Value = 3810.3713
$ws.Range("J17").Value = 3971.606
$ws.Range("L17").Value = 11914.818
$ws.Range("N17").Value = -12250.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2693862.5
$ws.Range("I2").Value = 3327437.5
$ws.Range("K2").Value = 3327437.5
$ws.Range("M2").Value = -3327324.5

$ws.Range("H24").Value = 26291.334
$ws.Range("J24").Value = 26291.334
$ws.Range("L24").Value = 26291.334
$ws.Range("N24").Value = -27039.334

$ws.Range("H100").Value = 26291.334
$ws.Range("J100").Value = 26291.334
$ws.Range("L100").Value = 26291.334
$ws.Range("N100").Value = -28455.334

$ws.Range("H116").Value = 2693862.5
$ws.Range("I116").Value = 3327437.5
$ws.Range("K116").Value = 3327437.5
$ws.Range("M116").Value = -3325143.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2693862.5
$ws.Range("I3").Value = 3327437.5
$ws.Range("K3").Value = 3327437.5
$ws.Range("M3").Value = -3327323.5

$ws.Range("H64").Value = 1188.375
$ws.Range("J64").Value = 1151.6666
$ws.Range("L64").Value = 1151.6666
$ws.Range("N64").Value = -1601.6666

$ws.Range("H67").Value = 1188.375
$ws.Range("J67").Value = 1151.6666
$ws.Range("L67").Value = 1151.6666
$ws.Range("N67").Value = -2711.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2394.3845
$ws.Range("I107").Value = 2394.3845
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2394.3845
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -474.3845000000001
$ws.Range("N107").ClearContents()

$ws.Range("H131").Value = 36383.87
$ws.Range("J131").Value = 36383.87
$ws.Range("L131").Value = 36383.87
$ws.Range("N131").Value = -46463.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 42714.5
$ws.Range("J5").Value = 144849.72
$ws.Range("L5").Value = 434549.16
$ws.Range("N5").Value = -434773.16

$ws.Range("H17").Value = 1251
$ws.Range("I17").Value = 3
$ws.Range("J17").Value = 1875
$ws.Range("K17").Value = 9
$ws.Range("L17").Value = 5625
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = -5963

$ws.Range("H25").Value = 449.64285
$ws.Range("I25").Value = 446.15384
$ws.Range("J25").Value = 495
$ws.Range("K25").Value = 1338.46152
$ws.Range("L25").Value = 1485
$ws.Range("M25").Value = -1169.46152
$ws.Range("N25").Value = -1823

$ws.Range("H30").Value = 449.64285
$ws.Range("I30").Value = 446.15384
$ws.Range("J30").Value = 495
$ws.Range("K30").Value = 1338.46152
$ws.Range("L30").Value = 1485
$ws.Range("M30").Value = -1236.46152
$ws.Range("N30").Value = -1689

$ws.Range("H32").Value = 100057860
$ws.Range("I32").Value = 68666.664
$ws.Range("J32").Value = 160051380
$ws.Range("K32").Value = 205999.992
$ws.Range("L32").Value = 480154140
$ws.Range("M32").Value = -205716.992
$ws.Range("N32").Value = -480154706

$ws.Range("H51").Value = 4855.364
$ws.Range("I51").Value = 568.1667
$ws.Range("K51").Value = 1704.5001
$ws.Range("M51").Value = -1244.5001

$ws.Range("I75").Value = 167.5
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 502.5
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = 495.5
$ws.Range("N75").Value = -4996

$ws.Range("I78").Value = 167.5
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 1507.5
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = 3484.5
$ws.Range("N78").Value = -18984

$ws.Range("H81").Value = 9454.929
$ws.Range("J81").Value = 9454.929
$ws.Range("L81").Value = 28364.787
$ws.Range("N81").Value = -30610.787

$ws.Range("H84").Value = 9454.929
$ws.Range("J84").Value = 9454.929
$ws.Range("L84").Value = 85094.361
$ws.Range("N84").Value = -96326.361

$ws.Range("H94").Value = 7838.6665
$ws.Range("I94").Value = 5137
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 15411
$ws.Range("L94").Value = 30000
$ws.Range("M94").Value = -14735
$ws.Range("N94").Value = -31352

$ws.Range("H114").Value = 879.1429000000001
$ws.Range("I114").Value = 1063
$ws.Range("J114").Value = 805.6
$ws.Range("K114").Value = 3189
$ws.Range("L114").Value = 2416.8
$ws.Range("M114").Value = 65
$ws.Range("N114").Value = -8924.799999999999

$ws.Range("H130").Value = 2176.6667
$ws.Range("I130").Value = 853.3333
$ws.Range("J130").Value = 3500
$ws.Range("K130").Value = 2559.9999
$ws.Range("L130").Value = 10500
$ws.Range("M130").Value = 2460.0001
$ws.Range("N130").Value = -20540

$ws.Range("H135").Value = 42714.5
$ws.Range("J135").Value = 144849.72
$ws.Range("L135").Value = 1303647.48
$ws.Range("N135").Value = -1308717.48

$ws.Range("H138").Value = 4408.1816
$ws.Range("I138").Value = 3998.5715
$ws.Range("K138").Value = 11995.7145
$ws.Range("M138").Value = -6855.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 27779412
$ws.Range("I113").Value = 33334874
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 33334874
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = -33332704
$ws.Range("N113").Value = -6440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3783.5
$ws.Range("I7").Value = 1922.2222
$ws.Range("K7").Value = 1922.2222
$ws.Range("M7").Value = -1810.2222

$ws.Range("H46").Value = 4412.25
$ws.Range("I46").Value = 2659.8
$ws.Range("K46").Value = 2659.8
$ws.Range("M46").Value = -2471.8

$ws.Range("H126").Value = 3783.5
$ws.Range("I126").Value = 1922.2222
$ws.Range("K126").Value = 5766.6666
$ws.Range("M126").Value = -3296.6666

$ws.Range("H136").Value = 68477.55
$ws.Range("I136").Value = 80111.69500000001
$ws.Range("J136").Value = 7980
$ws.Range("K136").Value = 240335.085
$ws.Range("L136").Value = 23940
$ws.Range("M136").Value = -237785.085
$ws.Range("N136").Value = -29040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 45457060
$ws.Range("J107").Value = 834.5
$ws.Range("L107").Value = 2503.5
$ws.Range("N107").Value = -6343.5

$ws.Range("H132").Value = 45952372
$ws.Range("I132").Value = 58825810
$ws.Range("K132").Value = 176477430
$ws.Range("M132").Value = -176474900
